$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Members")

# Fill in the (previously empty) "ID" column (C) = lowercase first name
$ids = @{
    2  = "albert"
    3  = "paul"
    4  = "isa"
    5  = "lea"
    6  = "elias"
    7  = "katell"
    8  = "dave"
    9  = "clara"
    10 = "bryan"
    11 = "fiona"
    12 = "alain"
}
foreach ($row in ($ids.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 3).Value = $ids[$row]
}

# New "Documentation" column (G)
$ws.Range("G1").Value = "Documentation"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G2").Value = "Some documentation about Albert`nMultiple lines`n* Item 1`n* Item 2"
$ws.Columns.Item(7).ColumnWidth = 55.5

# Match formatting: vertical-top alignment across the data block, wrap text on the doc cell
$ws.Range("A2:G12").VerticalAlignment = -4160
$ws.Range("G2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 90

# Window / selection state
$excel.ActiveWindow.Width = 37140
$excel.ActiveWindow.Height = 18570
$ws.Range("G3").Select() | Out-Null
